# Falguni Chouhan Invoice.xlsx - apply the "Add files via upload" edit.
#
# Summary of the change (from the OOXML diff):
#   - Sheet1 row 19: product changes from "Wires" to "Hard Drive"
#     (rate/amount recompute via existing VLOOKUP / multiplication formulas)
#   - Sheet1 row 22: previously-empty line item gets filled in with
#     SR.NO=10, "Scanner", qty=6 (rate/amount recompute via formulas)
#   - Totals (Total/SGST/CGST/Grand Total) recompute automatically because
#     they are formulas off of the row amounts.
#   - The active sheet / selection moves: Sheet1's selection moves to B19
#     and Sheet1 is no longer the active tab; Sheet2 becomes the active tab
#     with its previous scroll position reset.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Row 19: "Wires" -> "Hard Drive" -------------------------------------
# C19/E19 are formulas (VLOOKUP + multiplication); they recompute on their
# own once B19 changes because only a single precedent cell is touched.
$ws1.Range("B19").Value2 = "Hard Drive"

# --- Row 22: fill in the previously blank line item -----------------------
# Set D22 (quantity) before B22 (product) -- setting the VLOOKUP-driving
# cell (B22) last, after the other precedent (D22), ensures the dependent
# shared formula in E22 (=IFERROR(C22*D22,"")) recalculates correctly.
$ws1.Range("A22").Value2 = 10
$ws1.Range("D22").Value2 = 6
$ws1.Range("B22").Value2 = "Scanner"

# Nudge E22's shared formula so the cached value is refreshed (belt & braces
# in case the engine left a stale cached 0 behind after the B/D edits above).
$f22 = 'IFERROR(C22*D22,"")'
$ws1.Range("E22").Formula = "=" + $f22

# --- Selection / active sheet changes --------------------------------------
# Final state: Sheet1 selection at B19 (no longer the active tab),
# Sheet2 becomes the active tab.
$ws1.Range("B19").Select()
$ws2.Activate()
